$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Row 4: the first transaction's amount changes to 97, and the
# From/To names are swapped (Alan/Neady -> Neady/Alan).
$ws.Range("B4").Value = 97
$ws.Range("C4").Value = "Neady"
$ws.Range("D4").Value = "Alan"

# Row 5: this transaction is cleared out (amount removed) and its
# From/To names swapped back the other way (Neady/Alan -> Alan/Neady).
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = "Alan"
$ws.Range("D5").Value = "Neady"

# Rows 6 and 7: their transaction amounts are cleared out too.
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()

# New static figures recorded further down the "Alan View" block.
$ws.Range("J20").Value = 97
$ws.Range("J21").Value = 25
$ws.Range("J22").Value = 25
$ws.Range("J20:J22").Style = "Currency"

# Touch every formula so the whole sheet recalculates off the edited
# From/To names and amounts above (IF()-driven lookups on B4:D7 otherwise
# keep their previously cached results).
$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.HasFormula) {
        $cell.Formula = $cell.Formula
    }
}

# Update the saved selection to match the author's cursor position.
$ws.Range("B5").Select()
